$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 31.68831168831169

for ($r = 32; $r -le 51; $r++) {
    $ws.Cells.Item($r, 9).Value = $newValue
}
